$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-34 down to 9-35.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = "2021-09-14"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112022
$ws.Range("G8").Value = "Arveja Verde"
$ws.Range("H8").Value = "Perfection"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 35000
$ws.Range("L8").Value = 35000
$ws.Range("M8").Value = 35000
$ws.Range("N8").Value = "$/malla 25 kilos"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 1400
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
